# Applies the numeric value updates to the Leve profit-calculation sheets
# (currentAveragePrice* / LevePrice* / LeveProfit* columns H-N)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1467.85
$ws.Range("I19").Value = 2202.6667
$ws.Range("J19").Value = 365.625
$ws.Range("K19").Value = 2202.6667
$ws.Range("L19").Value = 365.625
$ws.Range("M19").Value = -2027.6667
$ws.Range("N19").Value = -715.625
$ws.Range("H62").Value = 10332.75
$ws.Range("I62").Value = 10399.3
$ws.Range("K62").Value = 10399.3
$ws.Range("M62").Value = -9775.299999999999
$ws.Range("H65").Value = 10332.75
$ws.Range("I65").Value = 10399.3
$ws.Range("K65").Value = 51996.5
$ws.Range("M65").Value = -48876.5
$ws.Range("H88").Value = 2715.0833
$ws.Range("I88").Value = 3000.25
$ws.Range("K88").Value = 3000.25
$ws.Range("M88").Value = -2594.25
$ws.Range("H91").Value = 2715.0833
$ws.Range("I91").Value = 3000.25
$ws.Range("K91").Value = 3000.25
$ws.Range("M91").Value = -1596.25
$ws.Range("H107").Value = 52250.53
$ws.Range("I107").Value = 41925.9
$ws.Range("K107").Value = 41925.9
$ws.Range("M107").Value = -40005.9
$ws.Range("H138").Value = 3919.2334
$ws.Range("I138").Value = 3345.4285
$ws.Range("J138").Value = 4421.3125
$ws.Range("K138").Value = 10036.2855
$ws.Range("L138").Value = 13263.9375
$ws.Range("M138").Value = -4896.2855
$ws.Range("N138").Value = -23543.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4274.7144
$ws.Range("I2").Value = 4351.9644
$ws.Range("K2").Value = 4351.9644
$ws.Range("M2").Value = -4238.9644
$ws.Range("H45").Value = 10733.167
$ws.Range("I45").Value = 12397.8
$ws.Range("K45").Value = 12397.8
$ws.Range("M45").Value = -12020.8
$ws.Range("H61").Value = 6672.9287
$ws.Range("I61").Value = 4301.654
$ws.Range("J61").Value = 37499.5
$ws.Range("K61").Value = 4301.654
$ws.Range("L61").Value = 37499.5
$ws.Range("M61").Value = -4089.654
$ws.Range("N61").Value = -37923.5
$ws.Range("H88").Value = 1960.375
$ws.Range("J88").Value = 1847.6428
$ws.Range("L88").Value = 1847.6428
$ws.Range("N88").Value = -2659.6428
$ws.Range("H91").Value = 1960.375
$ws.Range("J91").Value = 1847.6428
$ws.Range("L91").Value = 1847.6428
$ws.Range("N91").Value = -4655.6428
$ws.Range("H116").Value = 4274.7144
$ws.Range("I116").Value = 4351.9644
$ws.Range("K116").Value = 4351.9644
$ws.Range("M116").Value = -2057.9644
$ws.Range("H122").Value = 1799.3334
$ws.Range("I122").Value = 199.5
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 598.5
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = 1851.5
$ws.Range("N122").Value = -19897
$ws.Range("H132").Value = 3598.5715
$ws.Range("I132").Value = 3668.5
$ws.Range("K132").Value = 11005.5
$ws.Range("M132").Value = -8475.5
$ws.Range("H136").Value = 6672.9287
$ws.Range("I136").Value = 4301.654
$ws.Range("J136").Value = 37499.5
$ws.Range("K136").Value = 12904.962
$ws.Range("L136").Value = 112498.5
$ws.Range("M136").Value = -10354.962
$ws.Range("N136").Value = -117598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4274.7144
$ws.Range("I3").Value = 4351.9644
$ws.Range("K3").Value = 4351.9644
$ws.Range("M3").Value = -4237.9644
$ws.Range("H80").Value = 672.75
$ws.Range("I80").Value = 175.5
$ws.Range("K80").Value = 175.5
$ws.Range("M80").Value = 822.5
$ws.Range("H83").Value = 672.75
$ws.Range("I83").Value = 175.5
$ws.Range("K83").Value = 877.5
$ws.Range("M83").Value = 4114.5
$ws.Range("H87").Value = 100000
$ws.Range("I87").Value = 100000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 100000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -98752
$ws.Range("N87").Value = $null
$ws.Range("H90").Value = 100000
$ws.Range("I90").Value = 100000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 300000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -293760
$ws.Range("N90").Value = $null
$ws.Range("H105").Value = 3948.3142
$ws.Range("I105").Value = 3497.7407
$ws.Range("J105").Value = 5469
$ws.Range("K105").Value = 3497.7407
$ws.Range("L105").Value = 5469
$ws.Range("M105").Value = -1750.7407
$ws.Range("N105").Value = -8963
$ws.Range("H134").Value = 11765.393
$ws.Range("I134").Value = 5111.15
$ws.Range("K134").Value = 15333.45
$ws.Range("M134").Value = -12798.45

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = $null
$ws.Range("H31").Value = 3203.9048
$ws.Range("I31").Value = 2381.037
$ws.Range("J31").Value = 4685.067
$ws.Range("K31").Value = 2381.037
$ws.Range("L31").Value = 4685.067
$ws.Range("M31").Value = -2086.037
$ws.Range("N31").Value = -5275.067
$ws.Range("H34").Value = 3203.9048
$ws.Range("I34").Value = 2381.037
$ws.Range("J34").Value = 4685.067
$ws.Range("K34").Value = 2381.037
$ws.Range("L34").Value = 4685.067
$ws.Range("M34").Value = -2179.037
$ws.Range("N34").Value = -5089.067
$ws.Range("H58").Value = 5791.125
$ws.Range("I58").Value = 3794.4443
$ws.Range("J58").Value = 8358.286
$ws.Range("K58").Value = 3794.4443
$ws.Range("L58").Value = 8358.286
$ws.Range("M58").Value = -3591.4443
$ws.Range("N58").Value = -8764.286
$ws.Range("H87").Value = 64993.5
$ws.Range("J87").Value = 64993.5
$ws.Range("L87").Value = 64993.5
$ws.Range("N87").Value = -67365.5
$ws.Range("H90").Value = 64993.5
$ws.Range("J90").Value = 64993.5
$ws.Range("L90").Value = 194980.5
$ws.Range("N90").Value = -206836.5
$ws.Range("H99").Value = 8695.933999999999
$ws.Range("I99").Value = 4685.381
$ws.Range("J99").Value = 10855.462
$ws.Range("K99").Value = 4685.381
$ws.Range("L99").Value = 10855.462
$ws.Range("M99").Value = -3187.381
$ws.Range("N99").Value = -13851.462
$ws.Range("H105").Value = 707.4706
$ws.Range("I105").Value = 639.25
$ws.Range("K105").Value = 639.25
$ws.Range("M105").Value = 1107.75
$ws.Range("H126").Value = 8695.933999999999
$ws.Range("I126").Value = 4685.381
$ws.Range("J126").Value = 10855.462
$ws.Range("K126").Value = 14056.143
$ws.Range("L126").Value = 32566.386
$ws.Range("M126").Value = -11586.143
$ws.Range("N126").Value = -37506.386
$ws.Range("H134").Value = 10279.1
$ws.Range("I134").Value = 5559.4
$ws.Range("K134").Value = 16678.2
$ws.Range("M134").Value = -14143.2
$ws.Range("H136").Value = 5791.125
$ws.Range("I136").Value = 3794.4443
$ws.Range("J136").Value = 8358.286
$ws.Range("K136").Value = 11383.3329
$ws.Range("L136").Value = 25074.858
$ws.Range("M136").Value = -8833.332900000001
$ws.Range("N136").Value = -30174.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I132").Value = 2884.5386
$ws.Range("J132").Value = 8335764
$ws.Range("K132").Value = 25960.8474
$ws.Range("L132").Value = 75021876
$ws.Range("M132").Value = -23430.8474
$ws.Range("N132").Value = -75026936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 153325.14
$ws.Range("I113").Value = 258325.75
$ws.Range("J113").Value = 13324.333
$ws.Range("K113").Value = 258325.75
$ws.Range("L113").Value = 13324.333
$ws.Range("M113").Value = -256155.75
$ws.Range("N113").Value = -17664.333
$ws.Range("H116").Value = 120742
$ws.Range("J116").Value = 120742
$ws.Range("L116").Value = 120742
$ws.Range("N116").Value = -129920
$ws.Range("H132").Value = 7328.857
$ws.Range("I132").Value = 6431.154
$ws.Range("J132").Value = 18999
$ws.Range("K132").Value = 19293.462
$ws.Range("L132").Value = 56997
$ws.Range("M132").Value = -16763.462
$ws.Range("N132").Value = -62057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4246.8184
$ws.Range("I40").Value = 4191.3887
$ws.Range("J40").Value = 4496.25
$ws.Range("K40").Value = 4191.3887
$ws.Range("L40").Value = 4496.25
$ws.Range("M40").Value = -4055.3887
$ws.Range("N40").Value = -4768.25
$ws.Range("H93").Value = 4734.3667
$ws.Range("I93").Value = 5280.696
$ws.Range("J93").Value = 2939.2856
$ws.Range("K93").Value = 5280.696
$ws.Range("L93").Value = 2939.2856
$ws.Range("M93").Value = -4032.696
$ws.Range("N93").Value = -5435.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2357.1667
$ws.Range("I107").Value = 1684.2858
$ws.Range("J107").Value = 3299.2
$ws.Range("K107").Value = 5052.857400000001
$ws.Range("L107").Value = 9897.599999999999
$ws.Range("M107").Value = -3132.857400000001
$ws.Range("N107").Value = -13737.6
$ws.Range("H109").Value = 88855
$ws.Range("J109").Value = 89497.664
$ws.Range("L109").Value = 89497.664
$ws.Range("N109").Value = -92271.664
$ws.Range("H132").Value = 115377.17
$ws.Range("I132").Value = 135600.16
$ws.Range("K132").Value = 406800.48
$ws.Range("M132").Value = -404270.48
$ws.Range("H136").Value = 9093077
$ws.Range("I136").Value = 13045406
$ws.Range("K136").Value = 39136218
$ws.Range("M136").Value = -39133668
